$d = $word.ActiveDocument

# Locate the paragraph that holds "CEN 4010 Principles of Software
# Engineering, Summer 2021" - the new "Milestone 1" heading paragraph is
# inserted directly after it.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*CEN 4010 Principles of Software Engineering, Summer 2021*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $insertPos = $target.Range.End
    $insertRange = $d.Range($insertPos, $insertPos)

    # Insert the new "Milestone 1" heading paragraph right after it, followed
    # by a throw-away marker paragraph. InsertXML's last <w:p> always merges
    # into whatever paragraph already sits at the insertion point instead of
    # becoming a genuinely new one, so the marker paragraph is what actually
    # absorbs into (and thereby preserves) the pre-existing empty paragraph
    # that used to directly follow the CEN 4010 line.
    $xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:pPr>' +
            '<w:pStyle w:val="NormalWeb"/>' +
            '<w:spacing w:before="240" w:beforeAutospacing="0" w:after="240" w:afterAutospacing="0"/>' +
            '<w:jc w:val="center"/>' +
        '</w:pPr>' +
        '<w:r>' +
            '<w:rPr>' +
                '<w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>' +
                '<w:color w:val="000000"/>' +
                '<w:sz w:val="22"/>' +
                '<w:szCs w:val="22"/>' +
            '</w:rPr>' +
            '<w:t xml:space="preserve">     </w:t>' +
        '</w:r>' +
        '<w:r>' +
            '<w:rPr>' +
                '<w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>' +
                '<w:b/>' +
                '<w:bCs/>' +
                '<w:color w:val="000000"/>' +
                '<w:sz w:val="28"/>' +
                '<w:szCs w:val="28"/>' +
            '</w:rPr>' +
            '<w:t>Milestone 1: Team Project Proposal and Description</w:t>' +
        '</w:r>' +
        '</w:p>' +
        '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
            '<w:r><w:t>IRON_TMP_MARKER</w:t></w:r>' +
        '</w:p>'

    $insertRange.InsertXML($xml)

    # The marker paragraph now stands in for the original (pre-existing)
    # empty paragraph; strip its placeholder text back out so it reverts to
    # empty, without consuming its paragraph mark (which would merge it into
    # the next paragraph instead of leaving it standing alone).
    $marker = $null
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text -like "*IRON_TMP_MARKER*") {
            $marker = $p
            break
        }
    }
    if ($marker -ne $null) {
        $clearRange = $d.Range($marker.Range.Start, $marker.Range.End - 1)
        $clearRange.Delete()
    }
}
